$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (D:F) for Corequisites, Concurrent, Recommended.
# This shifts the existing "Terms Typically Offered" column (D) to G.
$ws.Columns("D:F").Insert()

# --- Header row ---
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
# G1 already holds the shifted "Terms Typically Offered" text from the column insert.

# --- Data rows 2-40 ---
# Row 2
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "F"

# Row 3
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "F, W, SP"

# Row 4
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "W"

# Row 5
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "TBD"

# Row 6
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "SP"

# Row 7
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("G7").Value = "F,W,SP,SU"

# Row 8
$ws.Range("C8").Value = "AGED 303; and junior standing."
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "F, SP"

# Row 9
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "F, W, SP"

# Row 10
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "F, W, SP"

# Row 11
$ws.Range("C11").Value = "Junior standing."
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "AGED 303."
$ws.Range("G11").Value = "W "

# Row 12
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "SP"

# Row 13
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("G13").Value = "W"

# Row 14
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "W, SP"

# Row 15
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "W"

# Row 16
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "F, W"

# Row 17
$ws.Range("C17").Value = "AGED 460."
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "W, SP"

# Row 18
$ws.Range("C18").Value = "AGED 461."
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "F, SP"

# Row 19
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "NA"
$ws.Range("G19").Value = "TBD"

# Row 20
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = "NA"
$ws.Range("G20").Value = "TBD"

# Row 21
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "TBD"

# Row 22
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "TBD"

# Row 23
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("F23").Value = "NA"
$ws.Range("G23").Value = "F, W, SP"

# Row 24
$ws.Range("D24").Value = "NA"
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "F"

# Row 25
$ws.Range("D25").Value = "NA"
$ws.Range("E25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "W, SP"

# Row 26
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = "NA"
$ws.Range("F26").Value = "NA"
$ws.Range("G26").Value = "SP"

# Row 27
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = "NA"
$ws.Range("F27").Value = "NA"
$ws.Range("G27").Value = "F"

# Row 28
$ws.Range("D28").Value = "NA"
$ws.Range("E28").Value = "NA"
$ws.Range("F28").Value = "NA"
$ws.Range("G28").Value = "TBD"

# Row 29
$ws.Range("C29").Value = "Admission to one of the following the Single Subject Agriculture Credential. the Agriculture Specialist Credential, or the Master of Agricultural Education. EDUC 412, EDUC 414, and EDUC 418."
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "AGED 303, and AGED 350."
$ws.Range("G29").Value = "F "

# Row 30
$ws.Range("C30").Value = "AEPS 120 or AEPS 230, and senior or graduate standing."
$ws.Range("D30").Value = "NA"
$ws.Range("E30").Value = "NA"
$ws.Range("F30").Value = "NA"
$ws.Range("G30").Value = "TBD"

# Row 31
$ws.Range("C31").Value = "AGED 438 and senior standing."
$ws.Range("D31").Value = "NA"
$ws.Range("E31").Value = "NA"
$ws.Range("F31").Value = "EDUC 410, EDUC 412, EDUC 414, EDUC 416 and EDUC 418."
$ws.Range("G31").Value = "TBD "

# Row 32
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = "NA"
$ws.Range("F32").Value = "NA"
$ws.Range("G32").Value = "SP"

# Row 33
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = "NA"
$ws.Range("F33").Value = "NA"
$ws.Range("G33").Value = "W"

# Row 34
$ws.Range("C34").Value = "AGED 438 or consent of instructor, enrollment in agriculture teaching credential program or MS degree in Agricultural Education, or current agriculture teacher; undergraduate biology course (BIO 111 or equivalent)."
$ws.Range("D34").Value = "NA"
$ws.Range("E34").Value = "NA"
$ws.Range("F34").Value = "NA"
$ws.Range("G34").Value = "TBD"

# Row 35
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = "NA"
$ws.Range("F35").Value = "NA"
$ws.Range("G35").Value = "F, W, SP"

# Row 36
$ws.Range("C36").Value = "AGED 524."
$ws.Range("D36").Value = "NA"
$ws.Range("E36").Value = "NA"
$ws.Range("F36").Value = "NA"
$ws.Range("G36").Value = "SU"

# Row 37
$ws.Range("D37").Value = "NA"
$ws.Range("E37").Value = "NA"
$ws.Range("F37").Value = "NA"
$ws.Range("G37").Value = "TBD"

# Row 38
$ws.Range("D38").Value = "NA"
$ws.Range("E38").Value = "NA"
$ws.Range("F38").Value = "NA"
$ws.Range("G38").Value = "TBD"

# Row 39
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = "NA"
$ws.Range("F39").Value = "NA"
$ws.Range("G39").Value = "F, W, SP"

# Row 40
$ws.Range("D40").Value = "NA"
$ws.Range("E40").Value = "NA"
$ws.Range("F40").Value = "NA"
$ws.Range("G40").Value = "TBD"
